$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for "Alcachofa" at "Macroferia Regional de Talca"
# needs to be inserted as row 126 (pushing existing rows 126-128 down to 127-129).
$ws.Rows.Item(126).Insert()

$ws.Range("A126").Value = 5
$ws.Range("B126").Value = "Macroferia Regional de Talca"
$ws.Range("C126").Value = "Maule"
$ws.Range("D126").Value = "2023-08-09"
$ws.Range("E126").Value = 7
$ws.Range("F126").Value = 100112013
$ws.Range("G126").Value = "Alcachofa"
$ws.Range("H126").Value = "Madrigal"
$ws.Range("I126").Value = "Primera"
$ws.Range("J126").Value = 300
$ws.Range("K126").Value = 13000
$ws.Range("L126").Value = 13000
$ws.Range("M126").Value = 13000
$ws.Range("N126").Value = "$/caja 40 unidades"
$ws.Range("O126").Value = "Provincia del Elquí"
$ws.Range("P126").Value = 325
$ws.Range("Q126").Value = 40
$ws.Range("R126").Value = "Hortaliza"
